$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timestamps = @(
  "2021-10-05 13:41:05.650305",
  "2021-10-05 13:41:05.650318",
  "2021-10-05 13:41:05.650322",
  "2021-10-05 13:41:05.650325",
  "2021-10-05 13:41:05.650328",
  "2021-10-05 13:41:05.650332",
  "2021-10-05 13:41:05.650335",
  "2021-10-05 13:41:05.650337",
  "2021-10-05 13:41:05.650341",
  "2021-10-05 13:41:05.650344",
  "2021-10-05 13:41:05.650347",
  "2021-10-05 13:41:05.650350",
  "2021-10-05 13:41:05.650353",
  "2021-10-05 13:41:05.650356",
  "2021-10-05 13:41:05.650359",
  "2021-10-05 13:41:05.650362",
  "2021-10-05 13:41:05.650365",
  "2021-10-05 13:41:05.650368",
  "2021-10-05 13:41:05.650371",
  "2021-10-05 13:41:05.650374",
  "2021-10-05 13:41:05.650377",
  "2021-10-05 13:41:05.650380",
  "2021-10-05 13:41:05.650383",
  "2021-10-05 13:41:05.650386",
  "2021-10-05 13:41:05.650389",
  "2021-10-05 13:41:05.650392",
  "2021-10-05 13:41:05.650396",
  "2021-10-05 13:41:05.650399",
  "2021-10-05 13:41:05.650402",
  "2021-10-05 13:41:05.650405",
  "2021-10-05 13:41:05.650408",
  "2021-10-05 13:41:05.650411",
  "2021-10-05 13:41:05.650414",
  "2021-10-05 13:41:05.650418",
  "2021-10-05 13:41:05.650421",
  "2021-10-05 13:41:05.650424",
  "2021-10-05 13:41:05.650427",
  "2021-10-05 13:41:05.650430",
  "2021-10-05 13:41:05.650433",
  "2021-10-05 13:41:05.650436",
  "2021-10-05 13:41:05.650439",
  "2021-10-05 13:41:05.650442"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
